$d = $word.ActiveDocument

# Disable smart-quote autocorrect so straight apostrophes survive inserts/replacements.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# 1. Remove the "Meta description" paragraph (2nd paragraph, right after the title).
#    It consists of an empty run, a bold "Meta description" run, and the description text run.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. The last paragraph currently holds the italic DALLE prompt text. Before it, insert a
#    new paragraph with an empty run followed by a bold run with the title text.
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)

$insertRange = $lastPara.Range.Duplicate
$insertRange.Collapse(1)
$insertRange.InsertBefore("Play Crystal Quest Frostlands Free - Exciting Slot Game`r")

# Re-fetch the newly created paragraph (it is now the second-to-last paragraph) and bold it.
$newPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$newPara.Range.Font.Bold = 1
$newPara.Range.Font.Italic = 0

# 3. Replace the text of the (now) last paragraph (the former DALLE prompt) with the new
#    meta description text, keeping its italic formatting.
$d.Content.Find.Execute(
    "Create a feature image fitting Crystal Quest Frostlands: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses DALLE, can you create a feature image for Crystal Quest Frostlands? The image should be in cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be standing in a snowy landscape, holding a crystal and fighting against a demonic yeti. The image should also include the game's logo, " + [char]34 + "Crystal Quest Frostlands," + [char]34 + " in bold letters and bright colors. The background should be icy with snowflakes falling, and the overall tone should be adventurous and exciting. Let's capture the thrill of the game with a visually stunning, attention-grabbing feature image!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Crystal Quest Frostlands features action-packed adventure, cascading reels and free spins. Play for free to enjoy the game's stunning design and big winning opportunities.",
    2
)
